# "break out stock.yaml completed"
#
# 1) Append 3 new rows (50-52) to the "day" sheet with the latest
#    day-timeframe stock-change data (TCS, MPHASIS, AUROPHARMA).
#    The bsecode column (D) is written as text, matching the newly
#    appended source rows.
# 2) On the "week" sheet, the bsecode column (D) for rows 67-71 was
#    previously stored as text; normalize it to a real number (same
#    values) like the rest of the column already is.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "day" sheet - append rows 50, 51, 52
# ---------------------------------------------------------------
$day = $wb.Worksheets.Item("day")

$newRows = @(
    @{ sr = 1; nsecode = "TCS";        name = "Tata Consultancy Services Limited"; bsecode = "532540"; per_chg = 0.16;  close = 3816.8;  volume = 1702163; tf = "day"; dt = "24/06/2024 11:34:51" },
    @{ sr = 2; nsecode = "MPHASIS";    name = "Mphasis Limited";                  bsecode = "526299"; per_chg = -1.46; close = 2394.2;  volume = 1356536; tf = "day"; dt = "24/06/2024 11:34:51" },
    @{ sr = 3; nsecode = "AUROPHARMA"; name = "Aurobindo Pharma Limited";          bsecode = "524804"; per_chg = -1.81; close = 1218.65; volume = 1623690; tf = "day"; dt = "24/06/2024 11:34:51" }
)

$r = 50
foreach ($row in $newRows) {
    $day.Cells.Item($r, 1).Value = $row.sr
    $day.Cells.Item($r, 2).Value = $row.nsecode
    $day.Cells.Item($r, 3).Value = $row.name

    # Keep bsecode as text (matches how this batch of rows was appended
    # upstream), not auto-converted to a number.
    $day.Cells.Item($r, 4).NumberFormat = "@"
    $day.Cells.Item($r, 4).Value = $row.bsecode

    $day.Cells.Item($r, 5).Value = $row.per_chg
    $day.Cells.Item($r, 6).Value = $row.close
    $day.Cells.Item($r, 7).Value = $row.volume
    $day.Cells.Item($r, 8).Value = $row.tf
    $day.Cells.Item($r, 9).Value = $row.dt

    $r = $r + 1
}

# ---------------------------------------------------------------
# 2) "week" sheet - fix bsecode type for rows 67-71 (text -> number)
# ---------------------------------------------------------------
$week = $wb.Worksheets.Item("week")

$bsecodes = @{ 67 = 532830; 68 = 532296; 69 = 532400; 70 = 532482; 71 = 500049 }

foreach ($rowNum in $bsecodes.Keys) {
    # Assigning a numeric literal (no NumberFormat change) flips the cell
    # from text to a genuine number without introducing a new style.
    $week.Cells.Item($rowNum, 4).Value = $bsecodes[$rowNum]
}
